# Apply the edits described by the diff to the workbook.
$wb = $excel.ActiveWorkbook

$wsUnits        = $wb.Worksheets.Item("units")
$wsIngredients  = $wb.Worksheets.Item("ingredients")
$wsStepTypes    = $wb.Worksheets.Item("step_types")
$wsSteps        = $wb.Worksheets.Item("steps")
$wsRecipes      = $wb.Worksheets.Item("recipes")
$wsRecipeSteps  = $wb.Worksheets.Item("recipe_steps")

# --- recipes sheet: update ingredient names to capitalized variants ---
$wsRecipes.Range("B2").Value = "White basmati rice"
$wsRecipes.Range("B3").Value = "Long grain brown rice"

# --- recipe_steps sheet: rename headers (simplify schema: step_id -> type, ingredient_id -> ingredient) ---
$wsRecipeSteps.Range("B1").Value = "type"
$wsRecipeSteps.Range("C1").Value = "ingredient"

# --- active tab moves from "steps" to "recipes" ---
$wsRecipes.Activate()

# --- selection changes on each sheet ---
$wsIngredients.Range("A26").Select()
$wsSteps.Range("C39").Select()
$wsRecipes.Range("E7").Select()
$wsRecipeSteps.Range("G7").Select()

# re-activate recipes as the final active sheet (its tabSelected should be "1")
$wsRecipes.Activate()
$wsRecipes.Range("E7").Select()
